# Applies the re-shuffling of species-observation rows described by the
# source diff: rows 3/4 swap their entire contents, rows 22/23 swap their
# entire contents, and rows 25/26/27 rotate their contents. Row numbers
# (and therefore unrelated formatting tied to the row) stay put; only the
# cell values move between rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 <-> Row 4 (full content swap)
$ws.Range("A3").Value2 = 130862972
$ws.Range("A4").Value2 = 130862980
$ws.Range("AC3").Value2 = "Flera fruktkroppar i en gammal relativt grov döende stående gran med full längd och 36 cm i brösthöjdsdiameter.."
$ws.Range("AC4").Value2 = "På en tydligt gammal gran med nedåthängande grenar."
$ws.Range("AH3").Value2 = "Granskog"
$ws.Range("AH4").ClearContents()
$ws.Range("AJ3").Value2 = "gran"
$ws.Range("AJ4").ClearContents()
$ws.Range("AK3").Value2 = "Picea abies"
$ws.Range("AK4").ClearContents()
$ws.Range("AM3").Value2 = "Trädstam på levande träd"
$ws.Range("AM4").ClearContents()
$ws.Range("AO3").Value2 = "Stem on living tree # Picea abies"
$ws.Range("AO4").ClearContents()
$ws.Range("B3").Value2 = 91828
$ws.Range("B4").Value2 = 79243
$ws.Range("E3").Value2 = 5432
$ws.Range("E4").Value2 = 6425
$ws.Range("F3").Value2 = "Granticka"
$ws.Range("F4").Value2 = "Garnlav"
$ws.Range("G3").Value2 = "Porodaedalea chrysoloma s.lat."
$ws.Range("G4").Value2 = "Alectoria sarmentosa"
$ws.Range("H3").ClearContents()
$ws.Range("H4").Value2 = "(Ach.) Ach."
$ws.Range("K3").Value2 = "teleomorf"
$ws.Range("K4").ClearContents()
$ws.Range("Q3").Value2 = 447167
$ws.Range("Q4").Value2 = 447218
$ws.Range("R3").Value2 = 7042999
$ws.Range("R4").Value2 = 7042948

# Row 22 <-> Row 23 (full content swap)
$ws.Range("A22").Value2 = 130865710
$ws.Range("A23").Value2 = 130865715
$ws.Range("AB22").Value2 = "14:59"
$ws.Range("AB23").Value2 = "15:22"
$ws.Range("B22").Value2 = 79243
$ws.Range("B23").Value2 = 91808
$ws.Range("E22").Value2 = 6425
$ws.Range("E23").Value2 = 1202
$ws.Range("F22").Value2 = "Garnlav"
$ws.Range("F23").Value2 = "Ullticka"
$ws.Range("G22").Value2 = "Alectoria sarmentosa"
$ws.Range("G23").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H22").Value2 = "(Ach.) Ach."
$ws.Range("H23").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q22").Value2 = 447169
$ws.Range("Q23").Value2 = 447136
$ws.Range("R22").Value2 = 7042897
$ws.Range("R23").Value2 = 7043068
$ws.Range("S22").Value2 = 6
$ws.Range("S23").Value2 = 8
$ws.Range("Z22").Value2 = "14:59"
$ws.Range("Z23").Value2 = "15:22"

# Rows 25/26/27 rotation (new25=old27, new26=old25, new27=old26)
$ws.Range("A25").Value2 = 130865712
$ws.Range("A26").Value2 = 130865713
$ws.Range("A27").Value2 = 130865703
$ws.Range("AB25").Value2 = "15:13"
$ws.Range("AB26").Value2 = "15:18"
$ws.Range("AB27").Value2 = "14:20"
$ws.Range("B25").Value2 = 79243
$ws.Range("B26").Value2 = 91804
$ws.Range("B27").Value2 = 89193
$ws.Range("E25").Value2 = 6425
$ws.Range("E26").Value2 = 1108
$ws.Range("E27").Value2 = 510
$ws.Range("F25").Value2 = "Garnlav"
$ws.Range("F26").Value2 = "Harticka"
$ws.Range("F27").Value2 = "Doftskinn"
$ws.Range("G25").Value2 = "Alectoria sarmentosa"
$ws.Range("G26").Value2 = "Pelloporus leporinus"
$ws.Range("G27").Value2 = "Cystostereum murrayi"
$ws.Range("H25").Value2 = "(Ach.) Ach."
$ws.Range("H26").Value2 = "(Fr.) Krieglst."
$ws.Range("H27").Value2 = "(Berk. & M.A.Curtis.) Pouzar"
$ws.Range("Q25").Value2 = 447165
$ws.Range("Q26").Value2 = 447144
$ws.Range("Q27").Value2 = 447410
$ws.Range("R25").Value2 = 7043032
$ws.Range("R26").Value2 = 7043043
$ws.Range("R27").Value2 = 7042768
$ws.Range("S25").Value2 = 10
$ws.Range("S26").Value2 = 13
$ws.Range("S27").Value2 = 8
$ws.Range("Z25").Value2 = "15:13"
$ws.Range("Z26").Value2 = "15:18"
$ws.Range("Z27").Value2 = "14:20"
